# The sheet has a blank spacer row (row 33, all cells styled s="4" with no
# values) sitting between the "NCAP_BND" info row (row 32) and the
# "Ograniczenie mocy calkowitej" bound-table header (old row 34). Deleting
# that spacer row shifts every row below it up by one (old 34->33, 35->34,
# ..., 40->39), which matches the target layout and shrinks the used range
# from B2:M40 down to B2:M39.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("33").Delete()

# The active selection moves to I35 after the edit.
[void]$ws.Range("I35").Select()
